# Apply the data/metadata changes described by the commit
# "modified data (because there were some logic problems) & added python
# requirements & modified report"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("données08")

# --- Corrected data values (rows whose A/C columns had logic problems) ---
$updates = @(
    @{ Row = 17; A = 8.7999999999999989;  C = 122 },
    @{ Row = 19; A = 45.33;               C = 118 },
    @{ Row = 21; A = 34.910000000000004;  C = 133 },
    @{ Row = 27; A = 46.82;               C = 131 },
    @{ Row = 33; A = 24.490000000000002;  C = 129 },
    @{ Row = 38; A = 38.340000000000003;  C = 119 },
    @{ Row = 44; A = 53.580000000000005;  C = 133 },
    @{ Row = 46; A = 29.56;               C = 133 },
    @{ Row = 48; A = 25.75;               C = 119 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 1).Value = $u.A
    $ws.Cells.Item($u.Row, 3).Value = $u.C
}

# --- Sheet view: selection moved to C27 ---
$ws.Range("C27").Select()

# --- Window size of the workbook view ---
$excel.Width = 25800
$excel.Height = 13200

# --- absPath metadata tweak (folder renamed to Fichiers_2025new) ---
$wb.Path = "D:\Jean\MIASHS\Projet\Fichiers_2025new\"
